$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the Vecino / Telefono / Casa info for row 64 (Control "82 / 83"),
# leaving the Control column (A) untouched. The vacated values used to be
# "Constanza Mejia/Jorge Sanz", "222-951-1424", "115 nº 1140"; they become
# "-" like the other blank entries in the sheet.
$ws.Range("B64").Value = "-"
$ws.Range("C64").Value = "-"
$ws.Range("D64").Value = "-"

# Reflect the selection left behind in the saved file after scrolling down.
$excel.ActiveWindow.ScrollRow = 52
$ws.Range("A68").Select()
